$d = $word.ActiveDocument

# Locate the "4a. Invalid QR Code" paragraph (the "4a" step heading under
# "Alternative Flows:"). We find it by content rather than a hard-coded
# index so the script is resilient to minor paragraph-count drift.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("4a. Invalid QR Code")) {
        $target = $p
        break
    }
}

$full = $target.Range
$insStart = $d.Range($full.Start, $full.Start)

# Insert a throwaway character immediately before "4a" so that the
# subsequent Find/Replace match fully spans (and so removes) both the
# <w:proofErr w:type="gramStart"/> and <w:proofErr w:type="gramEnd"/>
# markers that currently wrap just the "4a" run, while also merging the
# "4a" run and the ". " run into one run - matching the target edit:
#   <w:proofErr gramStart/><w:r>4a</w:r><w:proofErr gramEnd/><w:r>. </w:r>
#   -> <w:r>4a. </w:r>
$insStart.InsertBefore("X")

$found = $d.Content.Find.Execute("X4a. ", $true, $false, $false, $false, $false, $true, 1, $false, "4a. ", 2)
